# EPBDS-9313 Show error message if the same RuntimeContext properties
# are defined on the same level.
#
# Adds a new datatype (MyDatatype5) that declares two RuntimeContext
# properties of the same type on the same level, plus a new Spreadsheet
# (mySpr2), a test table for it, and a supporting datatype (MyDatatype6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Datatype MyDatatype5 (rows 18-20): two same-level context properties.
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "Datatype MyDatatype5"

$ws.Range("B19").Value = "String"
$ws.Range("C19").Value = "state:context.usState"

$ws.Range("B20").Value = "String"
$ws.Range("C20").Value = "state2:context.usState"

# ---------------------------------------------------------------------
# Spreadsheet mySpr2 (rows 22-26).
# ---------------------------------------------------------------------
$ws.Range("B22").Value = "Spreadsheet Double mySpr2(MyDatatype6 myVar,MyDatatype6 myVar1 )"
$ws.Range("B22:D22").Merge()

$ws.Range("B23").Value = "properties"
$ws.Range("C23").Value = "caProvinces"
$ws.Range("D23").Value = "MB"

$ws.Range("B24").Value = "Steps"
$ws.Range("C24").Value = "Values"
$ws.Range("D24").Value = "Values"
$ws.Range("C24:D24").Merge()

$ws.Range("B25").Value = "Step1"
$ws.Range("C25").Value = 200
$ws.Range("D25").Value = 100

$ws.Range("B26").Value = "RETURN"
$ws.Range("C26").Value = "'=`$Step1"
$ws.Range("D26").Value = "'=`$Step1"
$ws.Range("C26:D26").Merge()

$ws.Range("B22:D26").Font.Name = "Calibri"
$ws.Range("B22:D26").Font.Size = 11
$ws.Rows("22:26").RowHeight = 15

# ---------------------------------------------------------------------
# Test table for mySpr2 (rows 28-32).
# ---------------------------------------------------------------------
$ws.Range("B28").Value = "Test  mySpr2"

$ws.Range("B29").Value = "_context_.caProvince"
$ws.Range("C29").Value = "myVar.myProvince"
$ws.Range("D29").Value = "myVar1.myProvince"
$ws.Range("E29").Value = "_res_"

$ws.Range("B30").Value = "_context_.caProvince"
$ws.Range("C30").Value = "myVar.myProvince"
$ws.Range("D30").Value = "myVar.someField"
$ws.Range("E30").Value = "_res_"

$ws.Range("B31").Value = "BC"
$ws.Range("C31").Value = "AB"
$ws.Range("D31").Value = "MB"
$ws.Range("E31").Value = 100

$ws.Range("B32").Value = "BC"
$ws.Range("C32").Value = "BC"
$ws.Range("D32").Value = "AB"
$ws.Range("E32").Value = 200

# ---------------------------------------------------------------------
# Datatype MyDatatype6 (rows 34-36) - supports mySpr2's parameters.
# ---------------------------------------------------------------------
$ws.Range("B34").Value = "Datatype  MyDatatype6"

$ws.Range("B35").Value = "String "
$ws.Range("C35").Value = "myProvince:context.caProvince"

$ws.Range("B36").Value = "Double"
$ws.Range("C36").Value = "someField"

$ws.Range("I31").Select() | Out-Null
